# correction in sa algorithm and 746 logs
# Rewrite column C (Fitness) values with corrected/plateaued figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-252 (Generation 0-250), grouped by the plateau value
# they settle at after the sa-algorithm correction.
$ws.Range("C2:C17").Value = 7908
$ws.Range("C18:C32").Value = 7639
$ws.Range("C33:C74").Value = 7312
$ws.Range("C75:C83").Value = 7295
$ws.Range("C84:C252").Value = 7293
